# Update the time-range values in column C. The encrypted "Посыл / Заповедь"
# values in column B (rows 2-7) repeat in pairs, and only the 2:55-3:0 /
# 3:0-3:5 pair (rows 2-3) and the 18:55-19:0 / 19:0-19:5 pair (rows 6-7)
# change to new time ranges; the 10:55-11:0 / 11:0-11:5 pair (rows 4-5)
# is left untouched.
#
# New shared strings must be appended in the order: "22:40-22:45",
# "22:45-22:50" (rows 6-7) before "8:15-8:20", "8:20-8:25" (rows 2-3), so
# update row 6/7 first, then row 2/3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = "22:40-22:45"
$ws.Range("C7").Value = "22:45-22:50"

$ws.Range("C2").Value = "8:15-8:20"
$ws.Range("C3").Value = "8:20-8:25"

# Update the active selection to reflect the edit location.
$ws.Range("C8").Select()
